$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.837.92'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.51%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.661.98'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.92%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.56'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.13%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.515'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.55%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.36'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.81%  '

$ws.Range("E9").Value = '  -0.30%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0621'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.14%  '

$ws.Range("E11").Value = '  -1.30%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.899.18'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.81%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.666.34'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.83%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.13'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.66%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.548'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.16%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.96'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.96%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '247.84'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +5.23%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '27.771.68'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.36%  '

$ws.Range("E19").Value = '  -1.56%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.47'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.53%  '

$ws.Range("E21").Value = '  -0.10%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.47'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.66%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.33'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.35%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.04'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.09%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.64'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.80%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.17'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.21%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.24'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.14%  '

$ws.Range("E28").Value = '  +0.10%  '

$ws.Range("E29").Value = '  -0.89%  '

$ws.Range("E30").Value = '  +6.02%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0499'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.02%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.34'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.97%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.435.28'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -7.00%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.12'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.57%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.55'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.95%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.39'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.16%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.928'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.21%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.580'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.64%  '

$ws.Range("E39").Value = '  -1.82%  '

$ws.Range("E40").Value = '  -1.63%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '69.26'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.73%  '

$ws.Range("E42").Value = '  +0.03%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.42'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.34%  '

$ws.Range("E44").Value = '  -1.48%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.804.60'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.18%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.788'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.90%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.70'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.30%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '89.17'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.86%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0₆0109'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.90%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.101'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.25%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.80'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.34%  '

